$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3961645.2
$ws.Range("J17").Value = 4149995
$ws.Range("L17").Value = 12449985
$ws.Range("N17").Value = -12450321

$ws.Range("H86").Value = 2600.1667
$ws.Range("I86").Value = 2589.75
$ws.Range("J86").Value = 2621
$ws.Range("K86").Value = 2589.75
$ws.Range("L86").Value = 2621
$ws.Range("M86").Value = -1466.75
$ws.Range("N86").Value = -4867

$ws.Range("H89").Value = 2600.1667
$ws.Range("I89").Value = 2589.75
$ws.Range("J89").Value = 2621
$ws.Range("K89").Value = 12948.75
$ws.Range("L89").Value = 13105
$ws.Range("M89").Value = -7332.75
$ws.Range("N89").Value = -24337

$ws.Range("H100").Value = 57403.277
$ws.Range("J100").Value = 6166.6665
$ws.Range("L100").Value = 6166.6665
$ws.Range("N100").Value = -7248.6665

$ws.Range("H106").Value = 26699
$ws.Range("I106").Value = 10759.5
$ws.Range("K106").Value = 10759.5
$ws.Range("M106").Value = -10128.5

$ws.Range("H116").Value = 9150.638999999999
$ws.Range("I116").Value = 5658.64
$ws.Range("J116").Value = 17087
$ws.Range("K116").Value = 5658.64
$ws.Range("L116").Value = 17087
$ws.Range("M116").Value = -2216.64
$ws.Range("N116").Value = -23971

$ws.Range("H125").Value = 1399.6666

$ws.Range("H132").Value = 2419.8108
$ws.Range("I132").Value = 2309.647
$ws.Range("K132").Value = 6928.941
$ws.Range("M132").Value = -4398.941


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3847.6812
$ws.Range("I32").Value = 3847.6812
$ws.Range("K32").Value = 3847.6812
$ws.Range("M32").Value = -3560.6812

$ws.Range("H61").Value = 15973.379
$ws.Range("J61").Value = 21644
$ws.Range("L61").Value = 21644
$ws.Range("N61").Value = -22068

$ws.Range("H74").Value = 20776.895
$ws.Range("I74").Value = 23738.084
$ws.Range("K74").Value = 23738.084
$ws.Range("M74").Value = -22864.084

$ws.Range("H77").Value = 20776.895
$ws.Range("I77").Value = 23738.084
$ws.Range("K77").Value = 118690.42
$ws.Range("M77").Value = -114322.42

$ws.Range("H122").Value = 2990.0625
$ws.Range("I122").Value = 2401.1667
$ws.Range("K122").Value = 7203.500100000001
$ws.Range("M122").Value = -4753.500100000001

$ws.Range("H136").Value = 15973.379
$ws.Range("J136").Value = 21644
$ws.Range("L136").Value = 64932
$ws.Range("N136").Value = -70032


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 479639
$ws.Range("I86").Value = 1251812
$ws.Range("J86").Value = 4455.615
$ws.Range("K86").Value = 1251812
$ws.Range("L86").Value = 4455.615
$ws.Range("M86").Value = -1250689
$ws.Range("N86").Value = -6701.615

$ws.Range("H89").Value = 479639
$ws.Range("I89").Value = 1251812
$ws.Range("J89").Value = 4455.615
$ws.Range("K89").Value = 6259060
$ws.Range("L89").Value = 22278.075
$ws.Range("M89").Value = -6253444
$ws.Range("N89").Value = -33510.075


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3270.8333
$ws.Range("I31").Value = 1344.2858
$ws.Range("K31").Value = 1344.2858
$ws.Range("M31").Value = -1049.2858

$ws.Range("H34").Value = 3270.8333
$ws.Range("I34").Value = 1344.2858
$ws.Range("K34").Value = 1344.2858
$ws.Range("M34").Value = -1142.2858

$ws.Range("H58").Value = 3061.95
$ws.Range("I58").Value = 2159
$ws.Range("J58").Value = 4326.08
$ws.Range("K58").Value = 2159
$ws.Range("L58").Value = 4326.08
$ws.Range("M58").Value = -1956
$ws.Range("N58").Value = -4732.08

$ws.Range("H122").Value = 2419.0625
$ws.Range("I122").Value = 1256.75
$ws.Range("J122").Value = 5906
$ws.Range("K122").Value = 3770.25
$ws.Range("L122").Value = 17718
$ws.Range("M122").Value = -1320.25
$ws.Range("N122").Value = -22618

$ws.Range("H136").Value = 3061.95
$ws.Range("I136").Value = 2159
$ws.Range("J136").Value = 4326.08
$ws.Range("K136").Value = 6477
$ws.Range("L136").Value = 12978.24
$ws.Range("M136").Value = -3927
$ws.Range("N136").Value = -18078.24


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2410.4285
$ws.Range("I14").Value = 2410.4285
$ws.Range("K14").Value = 7231.2855
$ws.Range("M14").Value = -7058.2855

$ws.Range("H29").Value = 1270.5
$ws.Range("J29").Value = 2514.5
$ws.Range("L29").Value = 7543.5
$ws.Range("N29").Value = -8097.5

$ws.Range("H98").Value = 1410.7059
$ws.Range("I98").Value = 1297.2858
$ws.Range("J98").Value = 1490.1
$ws.Range("K98").Value = 3891.8574
$ws.Range("L98").Value = 4470.299999999999
$ws.Range("M98").Value = -2393.8574
$ws.Range("N98").Value = -7466.299999999999

$ws.Range("H131").Value = 3011.9395
$ws.Range("I131").Value = 1406.2727
$ws.Range("J131").Value = 3814.7727
$ws.Range("K131").Value = 4218.8181
$ws.Range("L131").Value = 11444.3181
$ws.Range("M131").Value = 821.1818999999996
$ws.Range("N131").Value = -21524.3181


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2180
$ws.Range("I80").Value = 2180
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2180
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1182
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 2180
$ws.Range("I83").Value = 2180
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 10900
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -5908
$ws.Range("N83").ClearContents()


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3110.75
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 3814.3333
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 3814.3333
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -4404.3333

$ws.Range("H27").Value = 3110.75
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 3814.3333
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 3814.3333
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -4028.3333

$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 6710.8
$ws.Range("I40").Value = 6710.8
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6710.8
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -6574.8
$ws.Range("N40").ClearContents()

$ws.Range("H82").Value = 2930.125
$ws.Range("I82").Value = 2625
$ws.Range("J82").Value = 3235.25
$ws.Range("K82").Value = 2625
$ws.Range("L82").Value = 3235.25
$ws.Range("M82").Value = -2264
$ws.Range("N82").Value = -3957.25

$ws.Range("H85").Value = 2930.125
$ws.Range("I85").Value = 2625
$ws.Range("J85").Value = 3235.25
$ws.Range("K85").Value = 2625
$ws.Range("L85").Value = 3235.25
$ws.Range("M85").Value = -1377
$ws.Range("N85").Value = -5731.25

$ws.Range("H98").Value = 32500
$ws.Range("J98").Value = 32500
$ws.Range("L98").Value = 32500
$ws.Range("N98").Value = -38490

$ws.Range("H99").Value = 31121.777
$ws.Range("J99").Value = 14000
$ws.Range("L99").Value = 14000
$ws.Range("N99").Value = -19990

$ws.Range("H122").Value = 4882.4375
$ws.Range("I122").Value = 4066.2727
$ws.Range("K122").Value = 12198.8181
$ws.Range("M122").Value = -9748.8181

$ws.Range("H132").Value = 4993.8823
$ws.Range("I132").Value = 4719.4194
$ws.Range("K132").Value = 14158.2582
$ws.Range("M132").Value = -11628.2582

$ws.Range("H136").Value = 4371.0684
$ws.Range("I136").Value = 4271.6943
$ws.Range("K136").Value = 12815.0829
$ws.Range("M136").Value = -10265.0829


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 19000
$ws.Range("J103").Value = 19000
$ws.Range("L103").Value = 19000
$ws.Range("N103").Value = -21344

$ws.Range("H122").Value = 3478.2727
$ws.Range("J122").Value = 19999
$ws.Range("L122").Value = 59997
$ws.Range("N122").Value = -64897

$ws.Range("H126").Value = 10138.767
$ws.Range("I126").Value = 6561.696
$ws.Range("J126").Value = 21892
$ws.Range("K126").Value = 19685.088
$ws.Range("L126").Value = 65676
$ws.Range("M126").Value = -17215.088
$ws.Range("N126").Value = -70616

$ws.Range("H132").Value = 158594.75
$ws.Range("I132").Value = 223143
$ws.Range("J132").Value = 32187.75
$ws.Range("K132").Value = 669429
$ws.Range("L132").Value = 96563.25
$ws.Range("M132").Value = -666899
$ws.Range("N132").Value = -101623.25

